$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as text so numeric-looking strings
# (e.g. '68.920.56', '1.00') are preserved exactly as in the source diff.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('B38').NumberFormat = '@'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('B39').NumberFormat = '@'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '68.920.56'
$ws.Range('E2').Value = '  +2.28%  '
$ws.Range('D3').Value = '2.526.32'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '595.90'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Value = '177.57'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.49%  '
$ws.Range('D9').Value = '2.524.82'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('E10').Value = '  +6.20%  '
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').Value = '4.99'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = '0.340'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '26.22'
$ws.Range('E14').Value = '  +2.77%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.979.83'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').Value = '68.812.98'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').Value = '2.535.96'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('D19').Value = '11.16'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').Value = '7.56'
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').Value = '355.40'
$ws.Range('E21').Value = '  +1.40%  '
$ws.Range('D22').Value = '4.12'
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '70.65'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').Value = '4.22'
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('D26').Value = '1.69'
$ws.Range('E26').Value = '  -5.61%  '
$ws.Range('D27').Value = '8.99'
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').Value = '523.43'
$ws.Range('E30').Value = '  +4.11%  '
$ws.Range('D31').Value = '0.0₃0893'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').Value = '7.79'
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '161.73'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').Value = '18.70'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').Value = '18.45'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('D40').Value = '1.78'
$ws.Range('E40').Value = '  +5.13%  '
$ws.Range('E41').Value = '  -1.12%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = '4.84'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = '0.327'
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = '153.11'
$ws.Range('E46').Value = '  +7.13%  '
$ws.Range('D47').Value = '3.58'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('D48').Value = '0.520'
$ws.Range('E48').Value = '  +2.17%  '
$ws.Range('D49').Value = '0.0₆0253'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D50').Value = '1.60'
$ws.Range('E50').Value = '  +1.98%  '
$ws.Range('D51').Value = '0.0740'
$ws.Range('E51').Value = '  -0.41%  '
